$wb = $excel.ActiveWorkbook

# Delete the first sheet ("version 1"), keeping only "version 2"
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("version 1").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheet to "Sheet1"
$ws = $wb.Worksheets.Item("version 2")
$ws.Name = "Sheet1"

# Select cell H9 on the remaining (now only) sheet
$ws.Activate()
[void]$ws.Range("H9").Select()
